$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.817
$ws.Range("D4").Value = -8.106999999999999
$ws.Range("B7").Value = 5.481
$ws.Range("A8").Value = -22.322
$ws.Range("A10").Value = -21.584
$ws.Range("D11").Value = -7.653
$ws.Range("A12").Value = -21.547
$ws.Range("B14").Value = 6.043
$ws.Range("D14").Value = -7.343999999999999
$ws.Range("B15").Value = 5.242000000000001
$ws.Range("A18").Value = -21.974
$ws.Range("B18").Value = 5.683
$ws.Range("D18").Value = -8.608000000000001
$ws.Range("D19").Value = -8.053999999999998
$ws.Range("B20").Value = 7.423
$ws.Range("D21").Value = -8.425000000000001
$ws.Range("A25").Value = -21.791
$ws.Range("D27").Value = -8.399000000000001
$ws.Range("B29").Value = 4.999
$ws.Range("B30").Value = 5.306
$ws.Range("B31").Value = 6.244
$ws.Range("D31").Value = -8.331
$ws.Range("B35").Value = 8.192
$ws.Range("A37").Value = -20.298
$ws.Range("D38").Value = -8.378
$ws.Range("B40").Value = 8.598000000000001
$ws.Range("D42").Value = -8.303999999999998
$ws.Range("B44").Value = 5.264
$ws.Range("D44").Value = -7.514999999999999
$ws.Range("D47").Value = -7.903
$ws.Range("B50").Value = 4.715000000000001
$ws.Range("B54").Value = 4.955
$ws.Range("A55").Value = -21.843
$ws.Range("D56").Value = -8.224
$ws.Range("D58").Value = -8.370999999999999
$ws.Range("D65").Value = -7.772
$ws.Range("A68").Value = -21.435
$ws.Range("B68").Value = 5.298
$ws.Range("D73").Value = -8.154
$ws.Range("B76").Value = 6.342999999999999
$ws.Range("A77").Value = -21.036
$ws.Range("A78").Value = -20.308
$ws.Range("A79").Value = -21.813
$ws.Range("A80").Value = -20.864
$ws.Range("A81").Value = -21.708
$ws.Range("A82").Value = -22.005
$ws.Range("A84").Value = -21.873
$ws.Range("B87").Value = 4.636
$ws.Range("B88").Value = 4.858
$ws.Range("D90").Value = -8.236999999999998
$ws.Range("B92").Value = 5.895999999999999
$ws.Range("D92").Value = -6.458999999999999
$ws.Range("D94").Value = -7.195
$ws.Range("D95").Value = -7.719000000000001
$ws.Range("B96").Value = 6.234000000000001
$ws.Range("B98").Value = 5.778
$ws.Range("A101").Value = -21.491
$ws.Range("B101").Value = 5.936999999999999
$ws.Range("D101").Value = -8.056000000000001
$ws.Range("A102").Value = -20.934
$ws.Range("B102").Value = 6.515000000000001
